# Barclaycard XLS format fix (#91)
#
# The source statement used two header labels that were longer / more
# verbose than the live bank export, and a stray duplicate label in the
# data row below the "Händlerdetails" column. Align the sheet text with
# the corrected format:
#   - "Name des Karteninhabers" -> "Karteninhaber"
#   - "Händlerdetails"          -> "Details"
#   - the data cell under that header ("Händler") now simply repeats the
#     (renamed) header text "Details", matching the real export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 13) label fixes.
$ws.Range("L13").Value = "Karteninhaber"
$ws.Range("O13").Value = "Details"

# Data row (row 14): drop the old "Händler" value in favour of the
# corrected header text.
$ws.Range("O14").Value = "Details"

# Move the active selection to the last data cell, matching the
# corrected workbook's saved cursor position.
$ws.Range("O14").Select() | Out-Null
